# Generate Report for Handoff
#
# The "ea69d92a-6c6a-4b68-af7d-dc59b83d9bb2.md" file (row 3 on every sheet)
# moves from "Handed back: in sync with en-US" to "Ready for handoff", with
# a refreshed handoff timestamp, and the per-language sheets record an
# error explaining that the handback file is stale.

$wb = $excel.ActiveWorkbook

$newStatus     = "Ready for handoff"
$overviewDate  = "2016-08-31 16:56:25"
$zhHandoffDate = "2016-08-31 16:56:21"
$deHandoffDate = "2016-08-31 16:56:25"
$errorDetail   = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e3e1cb1c8f04b14147297ac64b1c229ca14d48c2/e2e/ea69d92a-6c6a-4b68-af7d-dc59b83d9bb2.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7dfc26922a41574eef9b4c17e1995ae6630728f0/e2e/ea69d92a-6c6a-4b68-af7d-dc59b83d9bb2.md."

# ---- Overview sheet -------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus
$wsOverview.Range("G3").Value = $overviewDate

# ---- zh-cn sheet ------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = $newStatus
$wsZhCn.Range("H3").Value = $zhHandoffDate
$wsZhCn.Range("P3").Value = $errorDetail
$wsZhCn.Range("P1").EntireColumn.ColumnWidth = 39.166666666666664

# ---- de-de sheet ------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = $newStatus
$wsDeDe.Range("H3").Value = $deHandoffDate
$wsDeDe.Range("P3").Value = $errorDetail
$wsDeDe.Range("P1").EntireColumn.ColumnWidth = 39.166666666666664
